$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A92").Value = 1.705
$ws.Range("B92").Value = 1.518
$ws.Range("C92").Value = 2.348
